# Refactoring the codes and adding the docstrings
#
# The "Clusters Data" sheet lists, per cluster-count row, the set of
# resulting cluster ids. The ids used to overflow into columns B:X on row 2
# and B:K on row 3; after the refactor the same ids are re-wrapped so that
# row 2 only goes out to column P and row 3 goes out to column S.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Clusters Data")

# Final cluster-id layout (column -> value) for row 2 and row 3.
$row2 = [ordered]@{
    "B2" = "105"; "C2" = "98";  "D2" = "81";  "E2" = "13";  "F2" = "65";
    "G2" = "99";  "H2" = "86";  "I2" = "89";  "J2" = "103"; "K2" = "112";
    "L2" = "94";  "M2" = "97";  "N2" = "115"; "O2" = "106"; "P2" = "111"
}

$row3 = [ordered]@{
    "B3" = "109"; "C3" = "107"; "D3" = "113"; "E3" = "108"; "F3" = "66";
    "G3" = "100"; "H3" = "83";  "I3" = "40";  "J3" = "92";  "K3" = "114";
    "L3" = "120"; "M3" = "93";  "N3" = "192"; "O3" = "37";  "P3" = "75";
    "Q3" = "14";  "R3" = "74";  "S3" = "1"
}

# Drop the old overflow cells from row 2 (used to run through column X)
# that are no longer part of the data.
$ws.Range("Q2:X2").ClearContents()

# Write the text values for every cluster id cell, forcing text storage
# (these are numeric-looking strings, e.g. "105", which must stay text)
# by briefly applying the "@" text number format, then stripping the
# format back off so the cell keeps its original plain appearance.
foreach ($addr in $row2.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $row2[$addr]
    $cell.ClearFormats()
}

foreach ($addr in $row3.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $row3[$addr]
    $cell.ClearFormats()
}
